# Weekly fruit/vegetable price update:
# Insert one new record (row) at row 27, pushing the existing rows 27..79
# down to 28..80 (dimension grows from A1:R79 to A1:R80), and populate the
# newly inserted row with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 27..79 down by one row.
$ws.Rows.Item(27).Insert()

# Populate the new row 27 with this week's record.
$ws.Range("A27").Value = 10
$ws.Range("B27").Value = "Vega Modelo de Temuco"
$ws.Range("C27").Value = "La Araucanía"
$ws.Range("D27").Value = 44708
$ws.Range("E27").Value = 9
$ws.Range("F27").Value = 100112035
$ws.Range("G27").Value = "Bruselas (repollito)"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 20
$ws.Range("K27").Value = 30000
$ws.Range("L27").Value = 30000
$ws.Range("M27").Value = 30000
$ws.Range("N27").Value = "$/malla 10 kilos"
$ws.Range("O27").Value = "Región Metropolitana"
$ws.Range("P27").Value = 3000
$ws.Range("Q27").Value = 10
$ws.Range("R27").Value = "Hortaliza"
